$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Absent" (H) column values as part of forming the consolidated report.
$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
